$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.818.90"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.11%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.749.45"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.18%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.04%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'572.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.77%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'156.69"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.99%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.40%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E9").Value = "'  -4.25%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +0.23%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = "'Cardano"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = "'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = "'0.381"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -3.54%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = "'Toncoin"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = "'5.58"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -17.55%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'3.237.97"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.04%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'26.37"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -3.70%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'63.477.38"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.60%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.0000149"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -3.60%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'2.752.69"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.47%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'12.06"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.53%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  -3.33%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'353.75"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -2.98%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'6.71"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -5.30%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.19%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.532"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.69%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'65.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -3.03%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.170"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -1.45%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +0.23%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'8.37"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -3.21%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'0.0₃0897"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -2.16%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'1.92"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -5.06%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'6.98"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -3.39%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'169.21"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -3.24%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -5.61%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -3.29%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +0.26%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'4.84"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -1.60%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -1.42%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -2.53%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.974"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -3.97%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'6.16"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +4.67%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'4.12"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -4.56%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'327.27"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -4.57%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -1.26%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'21.35"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -3.74%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  -3.16%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'21.22"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -4.65%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value = "'Aave"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'134.81"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -3.08%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = "'VeChain"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'0.0252"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -3.27%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  -4.68%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.100"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.43%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'1.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.45%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'11.05"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +0.37%  "
$ws.Range("E51").Style = "Normal"
